$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020745528196649
$ws.Range("D2").Value = 1.025947322369441
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.019164524913674
$ws.Range("I2").Value = 1.030222605984031
$ws.Range("J2").Value = 1.025940902569356
$ws.Range("K2").Value = 1.028771403585333
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.022008583912765
$ws.Range("N2").Value = 1.012734953412092
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021496325061568
$ws.Range("D3").Value = 1.026479811865787
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.020551330542654
$ws.Range("I3").Value = 1.030347340266065
$ws.Range("J3").Value = 1.026330055782694
$ws.Range("K3").Value = 1.029112049596697
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.023199758582473
$ws.Range("N3").Value = 1.012863917257193
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.021982354899009
$ws.Range("D4").Value = 1.026824538590819
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.021448987420873
$ws.Range("I4").Value = 1.030426954316081
$ws.Range("J4").Value = 1.026581403887187
$ws.Range("K4").Value = 1.029331930714514
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.023970315791972
$ws.Range("N4").Value = 1.012947198624322
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02218673193511
$ws.Range("D5").Value = 1.026969501187289
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.021826437140035
$ws.Range("I5").Value = 1.030460160931685
$ws.Range("J5").Value = 1.026686959861647
$ws.Range("K5").Value = 1.029424238798991
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.02429420940559
$ws.Range("N5").Value = 1.012982169927577
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022221050642625
$ws.Range("D6").Value = 1.026993843321917
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.021889817072577
$ws.Range("I6").Value = 1.030465721026277
$ws.Range("J6").Value = 1.026704676660265
$ws.Range("K6").Value = 1.029439730091128
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.024348589781109
$ws.Range("N6").Value = 1.012988039400019
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.021985085599588
$ws.Range("D7").Value = 1.026826475434659
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.021454030625294
$ws.Range("I7").Value = 1.03042739905895
$ws.Range("J7").Value = 1.026582814767327
$ws.Range("K7").Value = 1.029333164651166
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.023974643862623
$ws.Range("N7").Value = 1.012947666070779
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020999218496318
$ws.Range("D8").Value = 1.026127243612881
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.019633141707946
$ws.Range("I8").Value = 1.030264987485298
$ws.Range("J8").Value = 1.026072513595983
$ws.Range("K8").Value = 1.028886637784334
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.022411192297806
$ws.Range("N8").Value = 1.012778571756249
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019263684901566
$ws.Range("D9").Value = 1.024896470843216
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.016426685124377
$ws.Range("I9").Value = 1.02997041139306
$ws.Range("J9").Value = 1.025169802030601
$ws.Range("K9").Value = 1.028095700352109
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.019654464098354
$ws.Range("N9").Value = 1.012479338319264
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018107866192767
$ws.Range("D10").Value = 1.024076949560036
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.014290342073182
$ws.Range("I10").Value = 1.02976841339557
$ws.Range("J10").Value = 1.024565680461182
$ws.Range("K10").Value = 1.027565695759985
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.017815343186879
$ws.Range("N10").Value = 1.01227901069907
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017607683260013
$ws.Range("D11").Value = 1.023722339706717
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.013365550065164
$ws.Range("I11").Value = 1.029679619499024
$ws.Range("J11").Value = 1.024303547729594
$ws.Range("K11").Value = 1.027335563431817
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.017018642705555
$ws.Range("N11").Value = 1.012192070488255
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017421938321139
$ws.Range("D12").Value = 1.023590660407152
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.013022076557896
$ws.Range("I12").Value = 1.029646438436448
$ws.Range("J12").Value = 1.024206098773263
$ws.Range("K12").Value = 1.027249987054767
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.016722657402938
$ws.Range("N12").Value = 1.012159747633088
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017461779175155
$ws.Range("D13").Value = 1.023618904297052
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.013095751243508
$ws.Range("I13").Value = 1.029653564888678
$ws.Range("J13").Value = 1.024227005570226
$ws.Range("K13").Value = 1.027268347765239
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.016786149801846
$ws.Range("N13").Value = 1.012166682321815
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01759232859567
$ws.Range("D14").Value = 1.023711454260107
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.013337157735221
$ws.Range("I14").Value = 1.029676880800548
$ws.Range("J14").Value = 1.024295494221471
$ws.Range("K14").Value = 1.027328491595429
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.016994177612051
$ws.Range("N14").Value = 1.01218939926929
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017672770433491
$ws.Range("D15").Value = 1.023768482498723
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.013485900832936
$ws.Range("I15").Value = 1.029691220137772
$ws.Range("J15").Value = 1.024337681598283
$ws.Range("K15").Value = 1.027365535628655
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.017122343002792
$ws.Range("N15").Value = 1.012203392042246
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018141067976763
$ws.Range("D16").Value = 1.024100489189102
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.014351722628965
$ws.Range("I16").Value = 1.029774278418679
$ws.Range("J16").Value = 1.024583065931237
$ws.Range("K16").Value = 1.027580955526523
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.017868209978912
$ws.Range("N16").Value = 1.012284776501441
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.01843489817674
$ws.Range("D17").Value = 1.024308815595911
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.014894896587076
$ws.Range("I17").Value = 1.029826023558333
$ws.Range("J17").Value = 1.024736843762408
$ws.Range("K17").Value = 1.027715912727934
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.018335977190931
$ws.Range("N17").Value = 1.012335774222017
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01860631267497
$ws.Range("D18").Value = 1.024430352732841
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.015211745951913
$ws.Range("I18").Value = 1.029856077519862
$ws.Range("J18").Value = 1.024826487234146
$ws.Range("K18").Value = 1.027794569487799
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.018608784409843
$ws.Range("N18").Value = 1.012365501305036
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018664765358918
$ws.Range("D19").Value = 1.024471797764902
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.015319787773413
$ws.Range("I19").Value = 1.029866303404158
$ws.Range("J19").Value = 1.024857044413479
$ws.Range("K19").Value = 1.027821378958869
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.018701799028713
$ws.Range("N19").Value = 1.012375634238772
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.018403370040639
$ws.Range("D20").Value = 1.02428646166588
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.014836616581352
$ws.Range("I20").Value = 1.029820485042124
$ws.Range("J20").Value = 1.024720350286388
$ws.Range("K20").Value = 1.027701439452849
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.018285793679337
$ws.Range("N20").Value = 1.012330304615497
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017553883786825
$ws.Range("D21").Value = 1.023684199527173
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.013266068590138
$ws.Range("I21").Value = 1.029670020337537
$ws.Range("J21").Value = 1.02427532825745
$ws.Range("K21").Value = 1.027310783352938
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.016932920155467
$ws.Range("N21").Value = 1.012182710504075
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017020040527721
$ws.Range("D22").Value = 1.023305757313445
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.012278805962439
$ws.Range("I22").Value = 1.029574265715734
$ws.Range("J22").Value = 1.023995055787664
$ws.Range("K22").Value = 1.02706461315248
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.016081993687135
$ws.Range("N22").Value = 1.012089742291297
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017303015703812
$ws.Range("D23").Value = 1.023506355039991
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.012802154446236
$ws.Range("I23").Value = 1.029625136075135
$ws.Range("J23").Value = 1.024143677840751
$ws.Range("K23").Value = 1.02719516446079
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.016533117302501
$ws.Range("N23").Value = 1.012139042537692
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018417616160833
$ws.Range("D24").Value = 1.024296562371037
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.014862950728638
$ws.Range("I24").Value = 1.029822988055041
$ws.Range("J24").Value = 1.024727803140602
$ws.Range("K24").Value = 1.027707979492922
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.018308469550828
$ws.Range("N24").Value = 1.012332776154043
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.019712154693428
$ws.Range("D25").Value = 1.025214485570642
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.017255390155309
$ws.Range("I25").Value = 1.030047557656502
$ws.Range("J25").Value = 1.025403585098715
$ws.Range("K25").Value = 1.02830065800082
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.020367363571796
$ws.Range("N25").Value = 1.012556846148521
